# Alteração nos rótulos da tabela para já transformar a primeira linha em
# cabeçalho automaticamente no Power BI.
# Prefix the year / interval labels in row 1 (B1:E1) with "Ano " for the
# plain-year sheets and with "Intervalo " for the interval sheet, leaving
# column A (the row/series label) untouched.

$wb = $excel.ActiveWorkbook

# Sheets whose header row uses the simple "Ano YYYY" pattern.
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet whose header row uses year intervals.
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("B1").Value = "Intervalo 2015"
$ws.Range("C1").Value = "Intervalo 2015-2030"
$ws.Range("D1").Value = "Intervalo 2031-2040"
$ws.Range("E1").Value = "Intervalo 2041-2050"

# Sheet with only a single year column in the header.
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Range("B1").Value = "Ano 2015"
